$d = $word.ActiveDocument
$members = $d | Get-Member
foreach ($m in $members) { Write-Host $m.Name }
